# edit.ps1
#
# Applies two changes described by the commit "Update icon for immediate
# pause":
#
#   1. The cached `datetimeFigureOut` date field (slide master + all 11
#      custom layouts) advances from 10/8/25 to 10/10/25.
#   2. On slide 3, inside the "Group 20" icon group, the "Graphic 22"
#      picture (the exclamation-mark / "immediate pause" icon) is
#      repositioned/resized: its top moves down slightly and its height
#      shrinks (width/left stay the same).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date placeholder text everywhere it appears:
#    the slide master and every custom (slide) layout.
# ---------------------------------------------------------------------

$oldDate = "10/8/25"
$newDate = "10/10/25"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isPlaceholder = $false
        try { $isPlaceholder = ($shp.Type -eq 14) } catch {}
        if (-not $isPlaceholder) { continue }

        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}
        if ($phType -ne 16) { continue }   # ppPlaceholderDate

        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Resize/reposition the "immediate pause" exclamation icon on slide 3.
#    NOTE: writes to Top/Left/Width/Height on a shape nested in a group
#    address the shape's own (child) coordinate space directly, in
#    points (1 pt = 12700 EMU), so we target the exact EMU values from
#    the OOXML rather than the group-translated ("visual") coordinates
#    that reading .Top/.Left would otherwise report.
# ---------------------------------------------------------------------

$slide3 = $p.Slides.Item(3)

$iconGroup = $null
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shp = $slide3.Shapes.Item($i)
    if ($shp.Name -eq "Group 20" -or $shp.Id -eq 21) {
        $iconGroup = $shp
        break
    }
}
if ($iconGroup -eq $null) {
    throw "Could not find 'Group 20' shape on slide 3"
}

$exclamationIcon = $null
for ($i = 1; $i -le $iconGroup.GroupItems.Count; $i++) {
    $item = $iconGroup.GroupItems.Item($i)
    if ($item.Name -eq "Graphic 22" -or $item.Id -eq 23) {
        $exclamationIcon = $item
        break
    }
}
if ($exclamationIcon -eq $null) {
    throw "Could not find 'Graphic 22' shape inside 'Group 20' on slide 3"
}

# Target OOXML: <a:off x="1135777" y="4843746"/><a:ext cx="914400" cy="820546"/>
# x and cx (cx=914400=72pt) are unchanged; only y and cy move. A tiny
# epsilon is added past the exact EMU->point quotient so the host's
# point->EMU re-quantization on save reliably lands on the target EMU
# value instead of one EMU short (observed empirically).
$exclamationIcon.Top = 381.397325
$exclamationIcon.Height = 64.6099215
